$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 29.51802670773275
$ws.Range("D2").Value = -0.4019732922672468
$ws.Range("E2").Value = 0.1615825276961694
$ws.Range("C3").Value = 29.789942127982
$ws.Range("D3").Value = -0.1900578720179986
$ws.Range("E3").Value = 0.03612199471600994
$ws.Range("C4").Value = 30.01810839594576
$ws.Range("D4").Value = -0.02189160405424317
$ws.Range("E4").Value = 0.0004792423280677562
$ws.Range("C5").Value = 30.33414104483045
$ws.Range("D5").Value = 0.1241410448304521
$ws.Range("E5").Value = 0.01541099901159631
$ws.Range("C6").Value = 30.44695441741961
$ws.Range("D6").Value = 0.2269544174196128
$ws.Range("E6").Value = 0.05150830758627582
$ws.Range("C7").Value = 30.26487025116176
$ws.Range("D7").Value = -0.115129748838239
$ws.Range("E7").Value = 0.013254859067556
$ws.Range("C8").Value = 30.30002772833511
$ws.Range("D8").Value = -0.1399722716648881
$ws.Range("E8").Value = 0.01959223683502924
$ws.Range("C9").Value = 30.48236632579598
$ws.Range("D9").Value = 0.002366325795982505
$ws.Range("E9").Value = 0.000005599497772732234
$ws.Range("C10").Value = 30.62898618471906
$ws.Range("D10").Value = -0.06101381528094407
$ws.Range("E10").Value = 0.003722685655137164
$ws.Range("C11").Value = 30.18588338749158
$ws.Range("D11").Value = -0.5641166125084247
$ws.Range("E11").Value = 0.3182275525079802
$ws.Range("C12").Value = 30.31240229444805
$ws.Range("D12").Value = -0.6275977055519526
$ws.Range("E12").Value = 0.3938788800140754
$ws.Range("C13").Value = 30.6685515859415
$ws.Range("D13").Value = -0.2814484140584987
$ws.Range("E13").Value = 0.07921320977604415
$ws.Range("C14").Value = 31.02168174071198
$ws.Range("D14").Value = 0.001681740711983792
$ws.Range("E14").Value = 0.000002828251822343753
$ws.Range("C15").Value = 31.34523840567913
$ws.Range("D15").Value = 0.225238405679125
$ws.Range("E15").Value = 0.0507323393928741
$ws.Range("C16").Value = 31.32977761628349
$ws.Range("D16").Value = 0.04977761628348887
$ws.Range("E16").Value = 0.002477811082866256
$ws.Range("C17").Value = 31.76660192069113
$ws.Range("D17").Value = 0.3866019206911346
$ws.Range("E17").Value = 0.1494610450820744
$ws.Range("C18").Value = 32.07149440212046
$ws.Range("D18").Value = 0.491494402120459
$ws.Range("E18").Value = 0.2415667473157475
$ws.Range("C19").Value = 32.04447878937641
$ws.Range("D19").Value = 0.3944787893764143
$ws.Range("E19").Value = 0.1556135152678815
$ws.Range("C20").Value = 31.70459556565881
$ws.Range("D20").Value = -0.1754044343411856
$ws.Range("E20").Value = 0.03076671558655128
$ws.Range("C21").Value = 32.3360231156244
$ws.Range("D21").Value = 0.05602311562439866
$ws.Range("E21").Value = 0.003138589484264741
$ws.Range("C22").Value = 32.29917956816637
$ws.Range("D22").Value = -0.1508204318336297
$ws.Range("E22").Value = 0.02274680265848254
$ws.Range("C23").Value = 33.31699750471096
$ws.Range("D23").Value = 0.4669975047109602
$ws.Range("E23").Value = 0.2180866694062633
$ws.Range("C24").Value = 32.97503252456637
$ws.Range("D24").Value = 0.07503252456637455
$ws.Range("E24").Value = 0.0056298797428036
$ws.Range("C25").Value = 33.14286805253617
$ws.Range("D25").Value = 0.04286805253616421
$ws.Range("E25").Value = 0.001837669928243335
$ws.Range("C26").Value = 33.09659682782137
$ws.Range("D26").Value = -0.3034031721786263
$ws.Range("E26").Value = 0.09205348488805318
$ws.Range("C27").Value = 33.733575576439
$ws.Range("D27").Value = 0.03357557643899867
$ws.Range("E27").Value = 0.001127319333211042
$ws.Range("C28").Value = 34.52395447622725
$ws.Range("D28").Value = 0.4239544762272516
$ws.Range("E28").Value = 0.1797373979131232
$ws.Range("C29").Value = 34.55653235092028
$ws.Range("D29").Value = 0.1565323509202798
$ws.Range("E29").Value = 0.02450237688462963
$ws.Range("C30").Value = 35.13333827817324
$ws.Range("D30").Value = 0.2333382781732425
$ws.Range("E30").Value = 0.05444675206085349
$ws.Range("C31").Value = 34.92158448717328
$ws.Range("D31").Value = -0.3784155128267201
$ws.Range("E31").Value = 0.1431983003479096
$ws.Range("C32").Value = 35.23280952070267
$ws.Range("D32").Value = -0.4671904792973365
$ws.Range("E32").Value = 0.218266943946075
$ws.Range("C33").Value = 36.00392293331047
$ws.Range("D33").Value = -0.2960770666895272
$ws.Range("E33").Value = 0.08766162941947471
$ws.Range("C34").Value = 36.38271371876629
$ws.Range("D34").Value = -0.4172862812337073
$ws.Range("E34").Value = 0.1741278405058566
$ws.Range("C35").Value = 37.41941067871166
$ws.Range("D35").Value = 0.1194106787116667
$ws.Range("E35").Value = 0.01425891019038088
$ws.Range("C36").Value = 37.87342213573044
$ws.Range("D36").Value = -0.02657786426955511
$ws.Range("E36").Value = 0.0007063828691308941
$ws.Range("C37").Value = 38.40034623284688
$ws.Range("D37").Value = -0.09965376715312146
$ws.Range("E37").Value = 0.009930873307808551
$ws.Range("C38").Value = 39.26824967697826
$ws.Range("D38").Value = 0.3682496769782588
$ws.Range("E38").Value = 0.135607824594592
$ws.Range("C39").Value = 39.99580818383873
$ws.Range("D39").Value = 0.5958081838387272
$ws.Range("E39").Value = 0.3549873919292026
$ws.Range("C40").Value = 40.15353084845763
$ws.Range("D40").Value = 0.2535308484576362
$ws.Range("E40").Value = 0.06427789111964888
$ws.Range("C41").Value = 39.99429156781264
$ws.Range("D41").Value = -0.1057084321873631
$ws.Range("E41").Value = 0.01117427263551035
$ws.Range("C42").Value = 41.1970682034644
$ws.Range("D42").Value = 0.5970682034644028
$ws.Range("E42").Value = 0.3564904395882095
$ws.Range("C43").Value = 41.43450599492073
$ws.Range("D43").Value = 0.5345059949207354
$ws.Range("E43").Value = 0.2856966586062052
$ws.Range("C44").Value = 41.18589790951645
$ws.Range("D44").Value = -0.01410209048355426
$ws.Range("E44").Value = 0.0001988689560063517
$ws.Range("C45").Value = 41.07101789536907
$ws.Range("D45").Value = -0.4289821046309328
$ws.Range("E45").Value = 0.1840256460935846
$ws.Range("C46").Value = 41.55498737686809
$ws.Range("D46").Value = -0.2450126231319061
$ws.Range("E46").Value = 0.06003118549397744
$ws.Range("C47").Value = 42.09600345691913
$ws.Range("D47").Value = -0.1039965430808749
$ws.Range("E47").Value = 0.01081528097277227
$ws.Range("C48").Value = 42.41385264992605
$ws.Range("D48").Value = -0.2861473500739535
$ws.Range("E48").Value = 0.08188030595434571
$ws.Range("C49").Value = 43.66967031332786
$ws.Range("D49").Value = -0.03032968667213964
$ws.Range("E49").Value = 0.0009198898936301651
$ws.Range("C50").Value = 43.9408469767513
$ws.Range("D50").Value = -0.2591530232487074
$ws.Range("E50").Value = 0.06716028945894506
$ws.Range("C51").Value = 45.44306438683569
$ws.Range("D51").Value = -0.1569356131643076
$ws.Range("E51").Value = 0.0246287866792572
$ws.Range("C52").Value = -0.488767684261834
$ws.Range("E52").Value = 4.612971651533983
$ws.Range("E53").Value = 0.09225943303067967
